$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.3635854824943578
$ws.Cells.Item(2, 3).Value = 0.07369815540170066
$ws.Cells.Item(2, 5).Value = 0.7564522531944817
$ws.Cells.Item(2, 6).Value = 2.128429429425935
$ws.Cells.Item(2, 7).Value = 0.3768400343502876
$ws.Cells.Item(2, 8).Value = 0.5583961173340981
$ws.Cells.Item(2, 10).Value = 0.03536192482111566
$ws.Cells.Item(2, 11).Value = 0.3609811030691787
$ws.Cells.Item(2, 15).Value = 1.812212060673531

$ws.Cells.Item(3, 2).Value = 0.3200487342874112
$ws.Cells.Item(3, 3).Value = 0.0685778009563478
$ws.Cells.Item(3, 5).Value = 0.7239441850905024
$ws.Cells.Item(3, 6).Value = 2.103880541756794
$ws.Cells.Item(3, 7).Value = 0.3833712804624447
$ws.Cells.Item(3, 8).Value = 0.5652742159223472
$ws.Cells.Item(3, 10).Value = 0.03591906208995432
$ws.Cells.Item(3, 11).Value = 0.3150504799169198
$ws.Cells.Item(3, 15).Value = 1.840229098700959

$ws.Cells.Item(4, 2).Value = 0.2932565559408147
$ws.Cells.Item(4, 3).Value = 0.06542564129104278
$ws.Cells.Item(4, 5).Value = 0.7043218742967952
$ws.Cells.Item(4, 6).Value = 2.090136089265798
$ws.Cells.Item(4, 7).Value = 0.3877223743540874
$ws.Cells.Item(4, 8).Value = 0.5697775894499983
$ws.Cells.Item(4, 10).Value = 0.03628662435452945
$ws.Cells.Item(4, 11).Value = 0.2867404327439544
$ws.Cells.Item(4, 15).Value = 1.858733383585957

$ws.Cells.Item(5, 2).Value = 0.2823240473414614
$ws.Cells.Item(5, 3).Value = 0.06413914788585373
$ws.Cells.Item(5, 5).Value = 0.6964105665113038
$ws.Cells.Item(5, 6).Value = 2.084868772607066
$ws.Cells.Item(5, 7).Value = 0.3895810040181651
$ws.Cells.Item(5, 8).Value = 0.5716831932477895
$ws.Cells.Item(5, 10).Value = 0.03644280568842451
$ws.Cells.Item(5, 11).Value = 0.2751773546019649
$ws.Cells.Item(5, 15).Value = 1.866600909162898

$ws.Cells.Item(6, 2).Value = 0.280507859403599
$ws.Cells.Item(6, 3).Value = 0.06392541133037355
$ws.Cells.Item(6, 5).Value = 0.6951020305556739
$ws.Cells.Item(6, 6).Value = 2.084014278016639
$ws.Cells.Item(6, 7).Value = 0.389894787777461
$ws.Cells.Item(6, 8).Value = 0.5720038714471158
$ws.Cells.Item(6, 10).Value = 0.03646912553822457
$ws.Cells.Item(6, 11).Value = 0.2732557362134571
$ws.Cells.Item(6, 15).Value = 1.86792703540781

$ws.Cells.Item(7, 2).Value = 0.2931091739730505
$ws.Cells.Item(7, 3).Value = 0.06540829899162759
$ws.Cells.Item(7, 5).Value = 0.7042148356910616
$ws.Cells.Item(7, 6).Value = 2.090063702009559
$ws.Cells.Item(7, 7).Value = 0.3877470944599075
$ws.Cells.Item(7, 8).Value = 0.5698030038898168
$ws.Cells.Item(7, 10).Value = 0.03628870478332935
$ws.Cells.Item(7, 11).Value = 0.2865845953202495
$ws.Cells.Item(7, 15).Value = 1.858838164860714

$ws.Cells.Item(8, 2).Value = 0.3485869732317894
$ws.Cells.Item(8, 3).Value = 0.07193443165657243
$ws.Cells.Item(8, 5).Value = 0.7451733581813329
$ws.Cells.Item(8, 6).Value = 2.119688927239508
$ws.Cells.Item(8, 7).Value = 0.3790211417601057
$ws.Cells.Item(8, 8).Value = 0.5607095235422648
$ws.Cells.Item(8, 10).Value = 0.03554873196435082
$ws.Cells.Item(8, 11).Value = 0.345167181142898
$ws.Cells.Item(8, 15).Value = 1.821601899333899

$ws.Cells.Item(9, 2).Value = 0.4568728782523692
$ws.Cells.Item(9, 3).Value = 0.08466257892695239
$ws.Cells.Item(9, 5).Value = 0.8281777008848934
$ws.Cells.Item(9, 6).Value = 2.188352435700722
$ws.Cells.Item(9, 7).Value = 0.364622926147284
$ws.Cells.Item(9, 8).Value = 0.5451009298165914
$ws.Cells.Item(9, 10).Value = 0.03430023309807417
$ws.Cells.Item(9, 11).Value = 0.4591585839508525
$ws.Cells.Item(9, 15).Value = 1.758929665603446

$ws.Cells.Item(10, 2).Value = 0.5360945000430206
$ws.Cells.Item(10, 3).Value = 0.09396682578915261
$ws.Cells.Item(10, 5).Value = 0.8908115162778643
$ws.Cells.Item(10, 6).Value = 2.245287502833946
$ws.Cells.Item(10, 7).Value = 0.3557103764668383
$ws.Cells.Item(10, 8).Value = 0.5349893646681423
$ws.Cells.Item(10, 10).Value = 0.03350704043295494
$ws.Cells.Item(10, 11).Value = 0.5423362223658046
$ws.Cells.Item(10, 15).Value = 1.719218806427847

$ws.Cells.Item(11, 2).Value = 0.5720561591146804
$ws.Cells.Item(11, 3).Value = 0.09818840476303592
$ws.Cells.Item(11, 5).Value = 0.91966710186054
$ws.Cells.Item(11, 6).Value = 2.272607857619732
$ws.Cells.Item(11, 7).Value = 0.3520200448822521
$ws.Cells.Item(11, 8).Value = 0.5306839080260062
$ws.Cells.Item(11, 10).Value = 0.03317326354218153
$ws.Cells.Item(11, 11).Value = 0.5800459034477399
$ws.Cells.Item(11, 15).Value = 1.702534347360896

$ws.Cells.Item(12, 2).Value = 0.5856622633842221
$ws.Cells.Item(12, 3).Value = 0.09978533241663001
$ws.Cells.Item(12, 5).Value = 0.9306463070849276
$ws.Cells.Item(12, 6).Value = 2.283158263628479
$ws.Cells.Item(12, 7).Value = 0.3506751703616544
$ws.Cells.Item(12, 8).Value = 0.5290958894107121
$ws.Cells.Item(12, 10).Value = 0.03305077170241155
$ws.Cells.Item(12, 11).Value = 0.5943064751035649
$ws.Cells.Item(12, 15).Value = 1.696415341544451

$ws.Cells.Item(13, 2).Value = 0.5827324822592459
$ws.Cells.Item(13, 3).Value = 0.09944148258199448
$ws.Cells.Item(13, 5).Value = 0.9282794162028125
$ws.Cells.Item(13, 6).Value = 2.280876928585087
$ws.Cells.Item(13, 7).Value = 0.3509624713086907
$ws.Cells.Item(13, 8).Value = 0.5294360128586248
$ws.Cells.Item(13, 10).Value = 0.03307697876138427
$ws.Cells.Item(13, 11).Value = 0.5912360757955071
$ws.Cells.Item(13, 15).Value = 1.697724317738448

$ws.Cells.Item(14, 2).Value = 0.5731757823151327
$ws.Cells.Item(14, 3).Value = 0.09831981954236824
$ws.Cells.Item(14, 5).Value = 0.9205693210772523
$ws.Cells.Item(14, 6).Value = 2.27347173864797
$ws.Cells.Item(14, 7).Value = 0.3519083462893278
$ws.Cells.Item(14, 8).Value = 0.5305524115153233
$ws.Cells.Item(14, 10).Value = 0.03316310778770415
$ws.Cells.Item(14, 11).Value = 0.5812195200915369
$ws.Cells.Item(14, 15).Value = 1.702026940847148

$ws.Cells.Item(15, 2).Value = 0.5673204686562769
$ws.Cells.Item(15, 3).Value = 0.09763254470043137
$ws.Cells.Item(15, 5).Value = 0.9158534660005841
$ws.Cells.Item(15, 6).Value = 2.268962530608661
$ws.Cells.Item(15, 7).Value = 0.3524945754246289
$ws.Cells.Item(15, 8).Value = 0.5312417561508695
$ws.Cells.Item(15, 10).Value = 0.03321637285631063
$ws.Cells.Item(15, 11).Value = 0.5750815570760039
$ws.Cells.Item(15, 15).Value = 1.704688358932259

$ws.Cells.Item(16, 2).Value = 0.5337427124259193
$ws.Cells.Item(16, 3).Value = 0.09369070506595278
$ws.Cells.Item(16, 5).Value = 0.8889330433959088
$ws.Cells.Item(16, 6).Value = 2.243530679883378
$ws.Cells.Item(16, 7).Value = 0.3559588929513282
$ws.Cells.Item(16, 8).Value = 0.5352766625940433
$ws.Cells.Item(16, 10).Value = 0.03352939892161899
$ws.Cells.Item(16, 11).Value = 0.539869154811214
$ws.Cells.Item(16, 15).Value = 1.720336996660947

$ws.Cells.Item(17, 2).Value = 0.5131236625482813
$ws.Cells.Item(17, 3).Value = 0.09126962227217916
$ws.Cells.Item(17, 5).Value = 0.8725112470175702
$ws.Cells.Item(17, 6).Value = 2.228293222749656
$ws.Cells.Item(17, 7).Value = 0.3581775382698282
$ws.Cells.Item(17, 8).Value = 0.5378273627093861
$ws.Cells.Item(17, 10).Value = 0.03372836735695728
$ws.Cells.Item(17, 11).Value = 0.5182340518552451
$ws.Cells.Item(17, 15).Value = 1.730290835536536

$ws.Cells.Item(18, 2).Value = 0.5012569613205073
$ws.Cells.Item(18, 3).Value = 0.08987605381265951
$ws.Cells.Item(18, 5).Value = 0.8631000451603938
$ws.Cells.Item(18, 6).Value = 2.219662693845336
$ws.Cells.Item(18, 7).Value = 0.3594878952716201
$ws.Cells.Item(18, 8).Value = 0.5393221608124463
$ws.Cells.Item(18, 10).Value = 0.03384535411046663
$ws.Cells.Item(18, 11).Value = 0.5057780909561984
$ws.Cells.Item(18, 15).Value = 1.736145898568921

$ws.Cells.Item(19, 2).Value = 0.4972378924952352
$ws.Cells.Item(19, 3).Value = 0.08940404284284398
$ws.Cells.Item(19, 5).Value = 0.8599194445673106
$ws.Cells.Item(19, 6).Value = 2.216763483693654
$ws.Cells.Item(19, 7).Value = 0.3599374353497353
$ws.Cells.Item(19, 8).Value = 0.5398330306099055
$ws.Cells.Item(19, 10).Value = 0.03388540060734435
$ws.Cells.Item(19, 11).Value = 0.5015586770550726
$ws.Cells.Item(19, 15).Value = 1.73815061181746

$ws.Cells.Item(20, 2).Value = 0.5153193431911518
$ws.Cells.Item(20, 3).Value = 0.09152745765995007
$ws.Cells.Item(20, 5).Value = 0.8742558376652312
$ws.Cells.Item(20, 6).Value = 2.229901439523587
$ws.Cells.Item(20, 7).Value = 0.3579378132044511
$ws.Cells.Item(20, 8).Value = 0.5375529689111502
$ws.Cells.Item(20, 10).Value = 0.03370692330116221
$ws.Cells.Item(20, 11).Value = 0.520538395710787
$ws.Cells.Item(20, 15).Value = 1.72921778639892

$ws.Cells.Item(21, 2).Value = 0.5759831430148381
$ws.Cells.Item(21, 3).Value = 0.0986493260213166
$ws.Cells.Item(21, 5).Value = 0.9228325456229243
$ws.Cells.Item(21, 6).Value = 2.275641260414773
$ws.Cells.Item(21, 7).Value = 0.3516290912335123
$ws.Cells.Item(21, 8).Value = 0.5302233481071994
$ws.Cells.Item(21, 10).Value = 0.03313770360023582
$ws.Cells.Item(21, 11).Value = 0.5841621547354805
$ws.Cells.Item(21, 15).Value = 1.70075774984879

$ws.Cells.Item(22, 2).Value = 0.6155613336171371
$ws.Cells.Item(22, 3).Value = 0.1032939595275622
$ws.Cells.Item(22, 5).Value = 0.9548846198292438
$ws.Cells.Item(22, 6).Value = 2.306728749246247
$ws.Cells.Item(22, 7).Value = 0.3478125227896243
$ws.Cells.Item(22, 8).Value = 0.5256799670268322
$ws.Cells.Item(22, 10).Value = 0.03278843654686803
$ws.Cells.Item(22, 11).Value = 0.6256312817508558
$ws.Cells.Item(22, 15).Value = 1.683317870648779

$ws.Cells.Item(23, 2).Value = 0.594444286302803
$ws.Cells.Item(23, 3).Value = 0.1008159766902992
$ws.Cells.Item(23, 5).Value = 0.9377499721815639
$ws.Cells.Item(23, 6).Value = 2.290027352121882
$ws.Cells.Item(23, 7).Value = 0.3498213734766651
$ws.Cells.Item(23, 8).Value = 0.5280822464790873
$ws.Cells.Item(23, 10).Value = 0.03297276113007364
$ws.Cells.Item(23, 11).Value = 0.603509009557115
$ws.Cells.Item(23, 15).Value = 1.692519499428684

$ws.Cells.Item(24, 2).Value = 0.5143267149586848
$ws.Cells.Item(24, 3).Value = 0.0914108954288082
$ws.Cells.Item(24, 5).Value = 0.873467015075164
$ws.Cells.Item(24, 6).Value = 2.229173960784294
$ws.Cells.Item(24, 7).Value = 0.3580460844350313
$ws.Cells.Item(24, 8).Value = 0.5376769339885783
$ws.Cells.Item(24, 10).Value = 0.03371661007025573
$ws.Cells.Item(24, 11).Value = 0.5194966567762549
$ws.Cells.Item(24, 15).Value = 1.729702499321476

$ws.Cells.Item(25, 2).Value = 0.4276357022296509
$ws.Cells.Item(25, 3).Value = 0.08122721640663144
$ws.Cells.Item(25, 5).Value = 0.8054339829545398
$ws.Cells.Item(25, 6).Value = 2.168641050033258
$ws.Cells.Item(25, 7).Value = 0.3682263324518829
$ws.Cells.Item(25, 8).Value = 0.5490853617166849
$ws.Cells.Item(25, 10).Value = 0.03461623732287400
$ws.Cells.Item(25, 11).Value = 0.4284191648848434
$ws.Cells.Item(25, 15).Value = 1.774773653076849
